$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.293.24"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "1.561.84"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'210.46"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'44.34"
$ws.Range("E8").Value = "  -4.29%  "
$ws.Range("D9").Value = "'23.71"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Value = "'0.243"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "'0.0586"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "'0.0893"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "1.786.31"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "1.548.51"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "28.286.99"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "'0.511"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "'61.00"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "'227.68"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'7.34"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'3.93"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("D24").Value = "'8.87"
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").Value = "'150.17"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'14.85"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").Value = "1.378.13"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'0.0473"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'0.779"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("D47").Value = "'62.07"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "'0.915"
$ws.Range("E48").Value = "  -6.33%  "
$ws.Range("D49").Value = "1.699.40"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'85.26"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  -2.21%  "
